$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Proximity sheet: append 4 new log rows (rows 5-8) for "Living Room Main Door"
# ---------------------------------------------------------------------------
$proximity = $wb.Worksheets.Item("Proximity")

# Pre-format the date/timestamp columns as Text so Excel's input parser does
# not silently convert the date-like / time-like strings into date or time
# serial numbers. The format is reset back to "Normal" style right after the
# values are written so the cells end up without any explicit style, exactly
# like the rest of the sheet.
$proximity.Range("A5:B8").NumberFormat = "@"

$proximity.Range("A5").Value = "2026-02-01"
$proximity.Range("B5").Value = "20:23:03"
$proximity.Range("C5").Value = "20:00"
$proximity.Range("D5").Value = "Living Room Main Door"
$proximity.Range("E5").Value = "ENTER"
$proximity.Range("F5").Value = "User ENTERED Living Room Main Door"

$proximity.Range("A6").Value = "2026-02-01"
$proximity.Range("B6").Value = "20:23:05"
$proximity.Range("C6").Value = "20:00"
$proximity.Range("D6").Value = "Living Room Main Door"
$proximity.Range("E6").Value = "EXIT"
$proximity.Range("F6").Value = "User EXITED Living Room Main Door"

$proximity.Range("A7").Value = "2026-02-01"
$proximity.Range("B7").Value = "20:23:07"
$proximity.Range("C7").Value = "20:00"
$proximity.Range("D7").Value = "Living Room Main Door"
$proximity.Range("E7").Value = "ENTER"
$proximity.Range("F7").Value = "User ENTERED Living Room Main Door"

$proximity.Range("A8").Value = "2026-02-01"
$proximity.Range("B8").Value = "20:23:09"
$proximity.Range("C8").Value = "20:00"
$proximity.Range("D8").Value = "Living Room Main Door"
$proximity.Range("E8").Value = "EXIT"
$proximity.Range("F8").Value = "User EXITED Living Room Main Door"

$proximity.Range("A5:B8").Style = "Normal"

# ---------------------------------------------------------------------------
# Camera sheet: append 2 new log rows (rows 3-4) for "Living Room Main Door"
# ---------------------------------------------------------------------------
$camera = $wb.Worksheets.Item("Camera")

$camera.Range("A3:B4").NumberFormat = "@"

$camera.Range("A3").Value = "2026-02-01"
$camera.Range("B3").Value = "20:23:04"
$camera.Range("C3").Value = "20:00"
$camera.Range("D3").Value = "Living Room Main Door"
$camera.Range("E3").Value = "Image Captured"
$camera.Range("F3").Value = "Active"

$camera.Range("A4").Value = "2026-02-01"
$camera.Range("B4").Value = "20:23:08"
$camera.Range("C4").Value = "20:00"
$camera.Range("D4").Value = "Living Room Main Door"
$camera.Range("E4").Value = "Image Captured"
$camera.Range("F4").Value = "Active"

$camera.Range("A3:B4").Style = "Normal"
